$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.639.49"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "3.123.02"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'589.05"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'146.82"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.117.53"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +15.79%  "
$ws.Range("E11").Value = "  +3.77%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  +4.97%  "
$ws.Range("D14").Value = "'36.23"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.651.18"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "63.561.45"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("D19").Value = "3.120.85"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").Value = "'14.39"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "'13.30"
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("D25").Value = "'82.47"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +9.97%  "
$ws.Range("D28").Value = "'2.72"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").Value = "'27.23"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("D34").Value = "0.0₃0862"
$ws.Range("E34").Value = "  +7.45%  "
$ws.Range("D35").Value = "'2.37"
$ws.Range("E35").Value = "  +8.18%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").Value = "'3.39"
$ws.Range("E37").Value = "  +12.15%  "
$ws.Range("D38").Value = "'6.10"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "'51.00"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "'450.10"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").Value = "'0.0373"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "2.899.00"
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("D44").Value = "'0.276"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("E46").Value = "  +4.57%  "
$ws.Range("D47").Value = "'35.69"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "'126.31"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'0.112"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  +3.08%  "
